# Updated cryptos list with refreshed price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = "'28.850.17"
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').Value = "'1.876.26"
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').Value = '  -0.56%  '
$ws.Range('D5').Value = "'324.96"
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('D7').Value = "'0.4598"
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = "'0.3879"
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').Value = "'0.07863"
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').Value = "'0.9842"
$ws.Range('E10').Value = '  -1.69%  '
$ws.Range('D11').Value = "'21.76"
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('D12').Value = "'1.880.57"
$ws.Range('E12').Value = '  -1.88%  '
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').Value = "'5.660"
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('D15').Value = "'0.06947"
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('D16').Value = "'88.22"
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D18').Value = "'0.000009946"
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('D19').Value = "'16.96"
$ws.Range('E19').Value = '  -1.29%  '
$ws.Range('D20').Value = "'1.002"
$ws.Range('D21').Value = "'28.862.45"
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').Value = "'5.258"
$ws.Range('E22').Value = '  -1.10%  '
$ws.Range('D24').Value = "'2.086"
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('D25').Value = "'155.74"
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').Value = "'19.26"
$ws.Range('E26').Value = '  -0.41%  '
$ws.Range('D27').Value = "'5.988"
$ws.Range('E27').Value = '  +2.51%  '
$ws.Range('D28').Value = "'1.929"
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').Value = "'117.30"
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('D30').Value = "'0.09327"
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').Value = "'0.9042"
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('D32').Value = "'5.262"
$ws.Range('E32').Value = '  -0.82%  '
$ws.Range('D33').Value = "'1.323"
$ws.Range('E33').Value = '  -0.98%  '
$ws.Range('D34').Value = "'3.265"
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').Value = "'1.190"
$ws.Range('E35').Value = '  +2.85%  '
$ws.Range('D36').Value = "'0.05766"
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('D37').Value = "'0.02070"
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').Value = "'7.689"
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('D40').Value = "'0.5656"
$ws.Range('E40').Value = '  +0.24%  '
$ws.Range('E41').Value = '  -1.01%  '
$ws.Range('D42').Value = "'9.671"
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').Value = "'2.247"
$ws.Range('E43').Value = '  +1.23%  '
$ws.Range('D44').Value = "'11.84"
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('D45').Value = "'0.5348"
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').Value = "'0.07034"
$ws.Range('E46').Value = '  -1.91%  '
$ws.Range('D47').Value = "'1.844"
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('D48').Value = "'113.08"
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').Value = "'2.507"
$ws.Range('D50').Value = "'1.059"
$ws.Range('E50').Value = '  -5.28%  '
$ws.Range('D51').Value = "'70.64"
$ws.Range('E51').Value = '  -0.08%  '
